# Updates the cryptos price/volume table (and swaps the TrustWalletToken /
# InternetComputer(DFINITY) rows) per the "Updated cryptos list" commit.
#
# D-column cells store plain-looking decimal numbers (e.g. "1.001") as TEXT
# in the workbook (no thousands grouping, dot used as a display separator).
# Excel's Range.Value setter auto-detects numeric-looking strings and would
# silently coerce them to real numbers (losing the exact text and precision),
# so we force NumberFormat = "@" (Text) on each D cell before assigning the
# string value. E-column percentage strings keep their padding spaces, which
# already prevents Excel from treating them as numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '27.485.60'
$ws.Range("E2").Value = '  -2.57%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.747.69'
$ws.Range("E3").Value = '  -3.00%  '

$ws.Range("E4").Value = '  +0.16%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '323.87'
$ws.Range("E5").Value = '  -0.22%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.001'
$ws.Range("E6").Value = '  +0.14%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4439'
$ws.Range("E7").Value = '  +3.22%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3605'
$ws.Range("E8").Value = '  -0.89%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07452'
$ws.Range("E9").Value = '  -1.61%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '42.04'
$ws.Range("E10").Value = '  -5.96%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.095'
$ws.Range("E11").Value = '  -2.68%  '

$ws.Range("E12").Value = '  +0.16%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '20.51'
$ws.Range("E13").Value = '  -5.82%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.996'
$ws.Range("E14").Value = '  -3.37%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.099'
$ws.Range("E15").Value = '  -3.91%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.752.21'
$ws.Range("E16").Value = '  -3.42%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '91.86'
$ws.Range("E17").Value = '  -1.32%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.00001056'
$ws.Range("E18").Value = '  -1.37%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06402'
$ws.Range("E19").Value = '  +0.48%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '1.001'
$ws.Range("E20").Value = '  +0.09%  '

$ws.Range("E21").Value = '  -3.17%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.830'
$ws.Range("E22").Value = '  -3.09%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '27.536.95'
$ws.Range("E23").Value = '  -2.41%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '11.12'
$ws.Range("E24").Value = '  -2.81%  '

$ws.Range("E25").Value = '  -3.32%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '161.76'
$ws.Range("E26").Value = '  +1.04%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '20.32'
$ws.Range("E27").Value = '  -0.64%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.952.89'
$ws.Range("E28").Value = '  -3.24%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.067'
$ws.Range("E29").Value = '  -7.86%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '123.99'
$ws.Range("E30").Value = '  -3.32%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.066'
$ws.Range("E31").Value = '  -10.00%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.653'
$ws.Range("E32").Value = '  +3.50%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.08973'

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.462'
$ws.Range("E34").Value = '  -7.42%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '11.89'
$ws.Range("E35").Value = '  -7.38%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.02296'
$ws.Range("E36").Value = '  -3.01%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.2076'
$ws.Range("E37").Value = '  -2.51%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.6307'
$ws.Range("E38").Value = '  -3.31%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.05956'
$ws.Range("E39").Value = '  -2.89%  '

$ws.Range("B40").Value = 'TrustWalletToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.203'
$ws.Range("E40").Value = '  +0.38%  '

$ws.Range("B41").Value = 'InternetComputer(DFINITY)'
$ws.Range("C41").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '4.896'
$ws.Range("E41").Value = '  -4.90%  '

$ws.Range("E42").Value = '  +0.13%  '

$ws.Range("E43").Value = '  -3.29%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '7.735'
$ws.Range("E44").Value = '  -3.04%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '13.15'
$ws.Range("E45").Value = '  -3.32%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.705'
$ws.Range("E46").Value = '  -0.26%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.5861'
$ws.Range("E47").Value = '  -2.92%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '120.72'
$ws.Range("E48").Value = '  -3.98%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.935'
$ws.Range("E49").Value = '  -3.08%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.144'
$ws.Range("E50").Value = '  -1.32%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.06851'
$ws.Range("E51").Value = '  -1.74%  '
